$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 309.53845
$ws.Range("I53").Value = 228.53334
$ws.Range("J53").Value = 420
$ws.Range("K53").Value = 228.53334
$ws.Range("L53").Value = 420
$ws.Range("M53").Value = 408.46666
$ws.Range("N53").Value = -1694

# Hunk 1: ALC!row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 904.9737
$ws.Range("J129").Value = 954.92535
$ws.Range("L129").Value = 2864.77605
$ws.Range("N129").Value = -12864.77605

# Hunk 2: ALC!row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 9344
$ws.Range("I135").Value = 9997.666999999999
$ws.Range("K135").Value = 89979.003
$ws.Range("M135").Value = -87444.003

# Hunk 3: ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 35715910
$ws.Range("J137").Value = 3102
$ws.Range("L137").Value = 9306
$ws.Range("N137").Value = -14406

# Hunk 4: ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2784.6667
$ws.Range("I61").Value = 2341.6
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2341.6
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2129.6
$ws.Range("N61").Value = -5424

# Hunk 5: ARM!row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 12994.182
$ws.Range("I74").Value = 2517
$ws.Range("J74").Value = 40933.332
$ws.Range("K74").Value = 2517
$ws.Range("L74").Value = 40933.332
$ws.Range("M74").Value = -1643
$ws.Range("N74").Value = -42681.332

# Hunk 6: ARM!row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 12994.182
$ws.Range("I77").Value = 2517
$ws.Range("J77").Value = 40933.332
$ws.Range("K77").Value = 12585
$ws.Range("L77").Value = 204666.66
$ws.Range("M77").Value = -8217
$ws.Range("N77").Value = -213402.66

# Hunk 7: ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2627.4348
$ws.Range("I132").Value = 2371.65
$ws.Range("K132").Value = 7114.950000000001
$ws.Range("M132").Value = -4584.950000000001

# Hunk 8: ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2784.6667
$ws.Range("I136").Value = 2341.6
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7024.799999999999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4474.799999999999
$ws.Range("N136").Value = -20100

# Hunk 9: ARM!row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 42500
$ws.Range("J139").Value = 42500
$ws.Range("L139").Value = 42500
$ws.Range("N139").Value = -52780

# Hunk 10: BSM!row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Hunk 11: BSM!row 133
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Hunk 12: BSM!row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4112.353
$ws.Range("I134").Value = 2537.2727
$ws.Range("K134").Value = 7611.8181
$ws.Range("M134").Value = -5076.8181

# Hunk 13: CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1965.3158
$ws.Range("I31").Value = 1378.8823
$ws.Range("J31").Value = 6950
$ws.Range("K31").Value = 1378.8823
$ws.Range("L31").Value = 6950
$ws.Range("M31").Value = -1083.8823
$ws.Range("N31").Value = -7540

# Hunk 14: CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1965.3158
$ws.Range("I34").Value = 1378.8823
$ws.Range("J34").Value = 6950
$ws.Range("K34").Value = 1378.8823
$ws.Range("L34").Value = 6950
$ws.Range("M34").Value = -1176.8823
$ws.Range("N34").Value = -7354

# Hunk 15: CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2551.3125
$ws.Range("I58").Value = 1214.125
$ws.Range("K58").Value = 1214.125
$ws.Range("M58").Value = -1011.125

# Hunk 16: CRP!row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2425.3845
$ws.Range("I122").Value = 1400.4762
$ws.Range("J122").Value = 6730
$ws.Range("K122").Value = 4201.4286
$ws.Range("L122").Value = 20190
$ws.Range("M122").Value = -1751.4286
$ws.Range("N122").Value = -25090

# Hunk 17: CRP!row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3883.1875
$ws.Range("I132").Value = 3032.2222
$ws.Range("K132").Value = 9096.6666
$ws.Range("M132").Value = -6566.6666

# Hunk 18: CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3227.1365
$ws.Range("I134").Value = 1801.7142
$ws.Range("J134").Value = 5721.625
$ws.Range("K134").Value = 5405.142599999999
$ws.Range("L134").Value = 17164.875
$ws.Range("M134").Value = -2870.142599999999
$ws.Range("N134").Value = -22234.875

# Hunk 19: CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2551.3125
$ws.Range("I136").Value = 1214.125
$ws.Range("K136").Value = 3642.375
$ws.Range("M136").Value = -1092.375

# Hunk 20: CUL!row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4147.143
$ws.Range("I134").Value = 2838.3333
$ws.Range("K134").Value = 8514.999899999999
$ws.Range("M134").Value = -3444.999899999999

# Hunk 21: CUL!row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1593
$ws.Range("I139").Value = 1593
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 4779
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 361
$ws.Range("N139").ClearContents()

# Hunk 22: GSM!row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 328.16666
$ws.Range("I107").Value = 328.16666
$ws.Range("K107").Value = 328.16666
$ws.Range("M107").Value = 1591.83334

# Hunk 23: GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3102.147
$ws.Range("I132").Value = 2520.7727
$ws.Range("K132").Value = 7562.3181
$ws.Range("M132").Value = -5032.3181

# Hunk 24: GSM!row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 14598.333
$ws.Range("J133").Value = 14598.333
$ws.Range("L133").Value = 14598.333
$ws.Range("N133").Value = -24718.333

# Hunk 25: LTW!row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3809.0908
$ws.Range("I40").Value = 1660
$ws.Range("K40").Value = 1660
$ws.Range("M40").Value = -1524

# Hunk 26: LTW!row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1487.3125
$ws.Range("I68").Value = 1399.75
$ws.Range("J68").Value = 1750
$ws.Range("K68").Value = 1399.75
$ws.Range("L68").Value = 1750
$ws.Range("M68").Value = -650.75
$ws.Range("N68").Value = -3248

# Hunk 27: LTW!row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1487.3125
$ws.Range("I71").Value = 1399.75
$ws.Range("J71").Value = 1750
$ws.Range("K71").Value = 6998.75
$ws.Range("L71").Value = 8750
$ws.Range("M71").Value = -3254.75
$ws.Range("N71").Value = -16238

# Hunk 28: LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3689.4783
$ws.Range("I132").Value = 2543.6
$ws.Range("K132").Value = 7630.799999999999
$ws.Range("M132").Value = -5100.799999999999

# Hunk 29: WVR!row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15754.223
$ws.Range("I62").Value = 18771.428
$ws.Range("J62").Value = 5194
$ws.Range("K62").Value = 18771.428
$ws.Range("L62").Value = 5194
$ws.Range("M62").Value = -18147.428
$ws.Range("N62").Value = -6442

# Hunk 30: WVR!row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 15754.223
$ws.Range("I65").Value = 18771.428
$ws.Range("J65").Value = 5194
$ws.Range("K65").Value = 93857.14
$ws.Range("L65").Value = 25970
$ws.Range("M65").Value = -90737.14
$ws.Range("N65").Value = -32210

# Hunk 31: WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 23815436
$ws.Range("I132").Value = 41673816
$ws.Range("K132").Value = 125021448
$ws.Range("M132").Value = -125018918

# Hunk 32: WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 23883484
$ws.Range("I136").Value = 37149396
$ws.Range("K136").Value = 111448188
$ws.Range("M136").Value = -111445638
